# #5: property boat&car done
# Rework the "汽車" (car) sheet (3rd worksheet) to carry the full record
# layout used by the other property sheets: proper header names in row 1
# and the trailing property_category/category/date/legislator_name/
# legislator_id/source_file/index columns (H:N) populated on the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Extend the header (bold/bordered) and data-row formatting out to column N
# by copying the look of the existing header/data cells before filling them.
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

$ws.Range("G3").Copy()
$ws.Range("H3:N3").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Header row (row 1): field names instead of the first record's data ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2 (TOYOTARAV4 record): append the extra tracking columns ---
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-09"
$ws.Range("K2").Value = "林佳龍"
$ws.Range("L2").Value = 1741
$ws.Range("M2").Value = "tmpf4911"
$ws.Range("N2").Value = 32

# --- Row 3 (TOYOTAPrevia record): append the extra tracking columns ---
$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2012-04-09"
$ws.Range("K3").Value = "林佳龍"
$ws.Range("L3").Value = 1741
$ws.Range("M3").Value = "tmpf4911"
$ws.Range("N3").Value = 33
